# Apply "[Fonds de solidarite] Add 2022-06-14 data" update:
# refresh nombre_aides (col C) and montant_total (col E) counters for a
# set of rows (regional breakdown by NAF section) in the published data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Row = 9;   C = 69578;  E = 191682734 },
    @{ Row = 15;  C = 57549;  E = 238162864 },
    @{ Row = 46;  C = 55752;  E = 174093404 },
    @{ Row = 70;  C = 10843;  E = 37621871 },
    @{ Row = 117; C = 19732;  E = 56667618 },
    @{ Row = 168; C = 285120; E = 1213662244 },
    @{ Row = 170; C = 367572; E = 2848183112 },
    @{ Row = 171; C = 115227; E = 448912934 },
    @{ Row = 174; C = 357373; E = 1020004006 },
    @{ Row = 175; C = 125692; E = 815698201 },
    @{ Row = 178; C = 75370;  E = 102786543 },
    @{ Row = 179; C = 235803; E = 813674902 },
    @{ Row = 180; C = 141531; E = 341241854 },
    @{ Row = 213; C = 3639;   E = 11320543 },
    @{ Row = 277; C = 19091;  E = 42542549 },
    @{ Row = 293; C = 61672;  E = 194886211 },
    @{ Row = 313; C = 220664; E = 1371206867 },
    @{ Row = 318; C = 58823;  E = 225842549 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
